$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# YOU_FINISHED (English, column C, row 36): new wording
$ws.Range("C36").Value = "<h4>The {{test_name}} is finished.</h4>Please clicken `"Conitnue`" to proceed."

# CONTINUE_MAIN_TEST (German, column B, row 35):
# straight opening quote before "Weiter" becomes a German low opening quote „
$ws.Range("B35").Value = "Nun geht der Test los.<br> Bitte klicken Sie auf „Weiter`", wenn Sie bereit sind."

# Update the visible selection to match the edited cell (B35), mirroring the
# author's last selection before saving.
$ws.Range("B35").Select() | Out-Null
